# Update "想去人数" (want-to-go count) values in column F by +1
# for the specific rows that changed, on both the "展览" and
# "全部类型" worksheets (rId1 / sheet1.xml and rId4 / sheet4.xml).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    4  = 1525
    7  = 11172
    11 = 1073
    13 = 12260
    14 = 12851
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
